$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price-tracker sheet (GitHub Actions scheduled
# update): new Price (column D) and 1h-Volume-change (column E) readings
# for each coin row, plus the RenderToken / Dai rows (27-28) trading places
# with each other and picking up new figures.
#
# The Price/Volume columns are stored as plain text (e.g. "64.812.83" /
# "  +1.77%  "), not numbers, so any cell whose new reading happens to look
# like a genuine number (e.g. "593.45") is explicitly formatted as Text
# first -- otherwise Excel would silently reinterpret it as a numeric value.

$ws.Range('D2').Value = '64.828.51'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').Value = '3.162.57'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '593.45'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.48'
$ws.Range('E6').Value = '  +4.91%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '3.159.00'
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.02'
$ws.Range('E11').Value = '  +4.71%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.467'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '38.78'
$ws.Range('E13').Value = '  +5.31%  '
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').Value = '3.686.90'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.34'
$ws.Range('E17').Value = '  +3.94%  '
$ws.Range('D18').Value = '64.469.79'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').Value = '3.159.26'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '476.28'
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.05'
$ws.Range('E21').Value = '  +5.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.757'
$ws.Range('E22').Value = '  +3.38%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.72'
$ws.Range('E23').Value = '  +4.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.57'
$ws.Range('E24').Value = '  +4.99%  '
$ws.Range('E25').Value = '  +10.80%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '82.53'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.05'
$ws.Range('E27').Value = '  +8.93%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.75'
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.45'
$ws.Range('E30').Value = '  +6.75%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.119'
$ws.Range('E33').Value = '  +7.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '27.90'
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('D35').Value = '0.0₃0886'
$ws.Range('E35').Value = '  +5.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.57'
$ws.Range('E36').Value = '  +7.95%  '
$ws.Range('E37').Value = '  +2.79%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.24'
$ws.Range('E38').Value = '  +3.90%  '
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '468.97'
$ws.Range('E40').Value = '  +6.86%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.39'
$ws.Range('E41').Value = '  +6.85%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '51.41'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.302'
$ws.Range('E43').Value = '  +9.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0378'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').Value = '2.905.20'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '38.49'
$ws.Range('E47').Value = '  +4.15%  '
$ws.Range('E48').Value = '  +3.72%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '26.07'
$ws.Range('E49').Value = '  +8.07%  '
$ws.Range('E50').Value = '  +6.53%  '
$ws.Range('E51').Value = '  +0.06%  '
